# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# F2: 284 -> 290
# F4: 170 -> 171
# F5: 15  -> 16

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 290
    $ws.Range("F4").Value = 171
    $ws.Range("F5").Value = 16
}
